$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Volume/number line: "Volume 32   Number  46" -> "...47"
$ws.Range("A8").Value = "Volume 32   Number  47"

# Report covering week line: dates shift by one week
$ws.Range("C9").Value = "Report Covering the Week  11/17/2025  Through  11/23/2025"


$ws.Range("C15").Value = 2
$ws.Range("C15").NumberFormat = "#,##0"

$ws.Range("F15").Value = 2
$ws.Range("F15").NumberFormat = "#,##0"

$ws.Range("I15").Value = 13

$ws.Range("K15").Value = 85.714285714285

$ws.Range("L15").Value = 160

$ws.Range("M15").Value = -53.571428571428

$ws.Range("N15").Value = -70.454545454545

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D16").Value = 3

$ws.Range("E16").Value = -100

$ws.Range("F16").Value = 9

$ws.Range("H16").Value = 80

$ws.Range("J16").Value = 113

$ws.Range("K16").Value = -20.353982300885

$ws.Range("L16").Value = -23.728813559322

$ws.Range("M16").Value = -59.276018099547

$ws.Range("N16").Value = -85.804416403785

$ws.Range("C17").Value = 7

$ws.Range("E17").Value = 16.666666666666

$ws.Range("F17").Value = 11

$ws.Range("H17").Value = -42.105263157894

$ws.Range("I17").Value = 146

$ws.Range("J17").Value = 228

$ws.Range("K17").Value = -35.964912280701

$ws.Range("L17").Value = -31.132075471698

$ws.Range("M17").Value = -7.594936708860

$ws.Range("N17").Value = -76.451612903225

$ws.Range("C18").Value = 6
$ws.Range("C18").NumberFormat = "#,##0"

$ws.Range("D18").Value = 4

$ws.Range("E18").Value = 50

$ws.Range("F18").Value = 9

$ws.Range("G18").Value = 12

$ws.Range("H18").Value = -25

$ws.Range("I18").Value = 69

$ws.Range("J18").Value = 91

$ws.Range("K18").Value = -24.175824175824

$ws.Range("L18").Value = 4.545454545454

$ws.Range("M18").Value = -29.591836734693

$ws.Range("N18").Value = -90.142857142857

$ws.Range("C19").Value = 4

$ws.Range("D19").Value = 7

$ws.Range("E19").Value = -42.857142857142

$ws.Range("G19").Value = 26

$ws.Range("H19").Value = 7.692307692307

$ws.Range("I19").Value = 327

$ws.Range("J19").Value = 317

$ws.Range("K19").Value = 3.154574132492

$ws.Range("L19").Value = 2.830188679245

$ws.Range("M19").Value = 149.618320610687

$ws.Range("N19").Value = 3.809523809523

$ws.Range("D20").Value = 3

$ws.Range("E20").Value = -33.333333333333

$ws.Range("G20").Value = 8

$ws.Range("H20").Value = 25

$ws.Range("I20").Value = 73

$ws.Range("J20").Value = 72

$ws.Range("K20").Value = 1.388888888888

$ws.Range("L20").Value = -29.807692307692

$ws.Range("M20").Value = 62.222222222222

$ws.Range("N20").Value = -74.914089347079

$ws.Range("C21").Value = 21

$ws.Range("D21").Value = 23

$ws.Range("E21").Value = -8.695652173913

$ws.Range("F21").Value = 69

$ws.Range("G21").Value = 70

$ws.Range("H21").Value = -1.428571428571

$ws.Range("I21").Value = 720

$ws.Range("J21").Value = 831

$ws.Range("K21").Value = -13.357400722021

$ws.Range("L21").Value = -12.832929782082

$ws.Range("M21").Value = 4.803493449781

$ws.Range("N21").Value = -72.891566265060

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("G22").Value = 2

$ws.Range("H22").Value = -50

$ws.Range("L22").Value = -50

$ws.Range("M22").Value = 25

$ws.Range("M23").Value = -50

$ws.Range("C24").Value = 10

$ws.Range("D24").Value = 20

$ws.Range("E24").Value = -50

$ws.Range("F24").Value = 56

$ws.Range("G24").Value = 60

$ws.Range("H24").Value = -6.666666666666

$ws.Range("I24").Value = 666

$ws.Range("J24").Value = 691

$ws.Range("K24").Value = -3.617945007235

$ws.Range("L24").Value = -4.310344827586

$ws.Range("M24").Value = 88.668555240793

$ws.Range("C25").Value = 2
$ws.Range("C25").NumberFormat = "#,##0"

$ws.Range("D25").Value = 5

$ws.Range("E25").Value = -60

$ws.Range("F25").Value = 7

$ws.Range("H25").Value = -63.157894736842

$ws.Range("I25").Value = 118

$ws.Range("J25").Value = 206

$ws.Range("K25").Value = -42.718446601941

$ws.Range("L25").Value = -31.791907514450

$ws.Range("C26").Value = 12

$ws.Range("D26").Value = 7

$ws.Range("E26").Value = 71.428571428571

$ws.Range("F26").Value = 31

$ws.Range("G26").Value = 27

$ws.Range("H26").Value = 14.814814814814

$ws.Range("I26").Value = 332

$ws.Range("J26").Value = 329

$ws.Range("K26").Value = 0.911854103343

$ws.Range("L26").Value = 13.310580204778

$ws.Range("M26").Value = -27.982646420824

$ws.Range("C27").Value = 2
$ws.Range("C27").NumberFormat = "#,##0"

$ws.Range("F27").Value = 2
$ws.Range("F27").NumberFormat = "#,##0"

$ws.Range("I27").Value = 17

$ws.Range("K27").Value = 41.666666666666

$ws.Range("L27").Value = 21.428571428571

$ws.Range("D28").Value = 2

$ws.Range("J28").Value = 27

$ws.Range("K28").Value = -3.703703703703

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
